# Update the NATMI LR-pair results with refreshed TPM-based values.
# The sheet grows from 12 data rows (3 sending clusters x 4 target clusters)
# to 16 data rows (4 sending clusters x 4 target clusters), adding
# "Resolving-Mac" as an additional sending cluster. Every numeric column is
# refreshed with the new TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Fgf1'
$ws.Cells.Item(2, 3).Value = 'Nrp1'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 1.190640666666667
$ws.Cells.Item(2, 8).Value = 3.571922
$ws.Cells.Item(2, 9).Value = 0.1136540143525372
$ws.Cells.Item(2, 10).Value = 0.1136540143525372
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 133.7780026666667
$ws.Cells.Item(2, 14).Value = 401.334008
$ws.Cells.Item(2, 15).Value = 0.50863533211804
$ws.Cells.Item(2, 16).Value = 0.5086353321180399
$ws.Cells.Item(2, 17).Value = 159.2815302803751
$ws.Cells.Item(2, 18).Value = 1433.533772523376
$ws.Cells.Item(2, 19).Value = 0.05780844733675127
$ws.Cells.Item(2, 20).Value = 0.05780844733675126

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Fgf1'
$ws.Cells.Item(3, 3).Value = 'Nrp1'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 1.190640666666667
$ws.Cells.Item(3, 8).Value = 3.571922
$ws.Cells.Item(3, 9).Value = 0.1136540143525372
$ws.Cells.Item(3, 10).Value = 0.1136540143525372
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 52.441971
$ws.Cells.Item(3, 14).Value = 157.325913
$ws.Cells.Item(3, 15).Value = 0.1993888292903622
$ws.Cells.Item(3, 16).Value = 0.1993888292903622
$ws.Cells.Item(3, 17).Value = 62.439543312754
$ws.Cells.Item(3, 18).Value = 561.9558898147861
$ws.Cells.Item(3, 19).Value = 0.02266134086590243
$ws.Cells.Item(3, 20).Value = 0.02266134086590242

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Fgf1'
$ws.Cells.Item(4, 3).Value = 'Nrp1'
$ws.Cells.Item(4, 4).Value = 'MuSCs'
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 1.190640666666667
$ws.Cells.Item(4, 8).Value = 3.571922
$ws.Cells.Item(4, 9).Value = 0.1136540143525372
$ws.Cells.Item(4, 10).Value = 0.1136540143525372
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 21.197691
$ws.Cells.Item(4, 14).Value = 63.593073
$ws.Cells.Item(4, 15).Value = 0.08059542216956049
$ws.Cells.Item(4, 16).Value = 0.08059542216956046
$ws.Cells.Item(4, 17).Value = 25.238832944034
$ws.Cells.Item(4, 18).Value = 227.149496496306
$ws.Cells.Item(4, 19).Value = 0.009159993268008027
$ws.Cells.Item(4, 20).Value = 0.009159993268008024

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Fgf1'
$ws.Cells.Item(5, 3).Value = 'Nrp1'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 1.190640666666667
$ws.Cells.Item(5, 8).Value = 3.571922
$ws.Cells.Item(5, 9).Value = 0.1136540143525372
$ws.Cells.Item(5, 10).Value = 0.1136540143525372
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 55.59592133333333
$ws.Cells.Item(5, 14).Value = 166.787764
$ws.Cells.Item(5, 15).Value = 0.2113804164220374
$ws.Cells.Item(5, 16).Value = 0.2113804164220373
$ws.Cells.Item(5, 17).Value = 66.19476484026755
$ws.Cells.Item(5, 18).Value = 595.752883562408
$ws.Cells.Item(5, 19).Value = 0.02402423288187553
$ws.Cells.Item(5, 20).Value = 0.02402423288187553

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Fgf1'
$ws.Cells.Item(6, 3).Value = 'Nrp1'
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 5.168173666666667
$ws.Cells.Item(6, 8).Value = 15.504521
$ws.Cells.Item(6, 9).Value = 0.4933341355895272
$ws.Cells.Item(6, 10).Value = 0.4933341355895272
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 133.7780026666667
$ws.Cells.Item(6, 14).Value = 401.334008
$ws.Cells.Item(6, 15).Value = 0.50863533211804
$ws.Cells.Item(6, 16).Value = 0.5086353321180399
$ws.Cells.Item(6, 17).Value = 691.3879505611299
$ws.Cells.Item(6, 18).Value = 6222.491555050169
$ws.Cells.Item(6, 19).Value = 0.2509271719007454
$ws.Cells.Item(6, 20).Value = 0.2509271719007453

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Fgf1'
$ws.Cells.Item(7, 3).Value = 'Nrp1'
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 5.168173666666667
$ws.Cells.Item(7, 8).Value = 15.504521
$ws.Cells.Item(7, 9).Value = 0.4933341355895272
$ws.Cells.Item(7, 10).Value = 0.4933341355895272
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 52.441971
$ws.Cells.Item(7, 14).Value = 157.325913
$ws.Cells.Item(7, 15).Value = 0.1993888292903622
$ws.Cells.Item(7, 16).Value = 0.1993888292903622
$ws.Cells.Item(7, 17).Value = 271.029213550297
$ws.Cells.Item(7, 18).Value = 2439.262921952673
$ws.Cells.Item(7, 19).Value = 0.09836531574416865
$ws.Cells.Item(7, 20).Value = 0.09836531574416864

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Fgf1'
$ws.Cells.Item(8, 3).Value = 'Nrp1'
$ws.Cells.Item(8, 4).Value = 'MuSCs'
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 5.168173666666667
$ws.Cells.Item(8, 8).Value = 15.504521
$ws.Cells.Item(8, 9).Value = 0.4933341355895272
$ws.Cells.Item(8, 10).Value = 0.4933341355895272
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 21.197691
$ws.Cells.Item(8, 14).Value = 63.593073
$ws.Cells.Item(8, 15).Value = 0.08059542216956049
$ws.Cells.Item(8, 16).Value = 0.08059542216956046
$ws.Cells.Item(8, 17).Value = 109.553348420337
$ws.Cells.Item(8, 18).Value = 985.9801357830331
$ws.Cells.Item(8, 19).Value = 0.03976047292849314
$ws.Cells.Item(8, 20).Value = 0.03976047292849313

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Fgf1'
$ws.Cells.Item(9, 3).Value = 'Nrp1'
$ws.Cells.Item(9, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 5.168173666666667
$ws.Cells.Item(9, 8).Value = 15.504521
$ws.Cells.Item(9, 9).Value = 0.4933341355895272
$ws.Cells.Item(9, 10).Value = 0.4933341355895272
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 55.59592133333333
$ws.Cells.Item(9, 14).Value = 166.787764
$ws.Cells.Item(9, 15).Value = 0.2113804164220374
$ws.Cells.Item(9, 16).Value = 0.2113804164220373
$ws.Cells.Item(9, 17).Value = 287.3293766090049
$ws.Cells.Item(9, 18).Value = 2585.964389481044
$ws.Cells.Item(9, 19).Value = 0.1042811750161201
$ws.Cells.Item(9, 20).Value = 0.1042811750161201

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = 'MuSCs'
$ws.Cells.Item(10, 2).Value = 'Fgf1'
$ws.Cells.Item(10, 3).Value = 'Nrp1'
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 4.037194
$ws.Cells.Item(10, 8).Value = 12.111582
$ws.Cells.Item(10, 9).Value = 0.3853751326204581
$ws.Cells.Item(10, 10).Value = 0.3853751326204581
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 133.7780026666667
$ws.Cells.Item(10, 14).Value = 401.334008
$ws.Cells.Item(10, 15).Value = 0.50863533211804
$ws.Cells.Item(10, 16).Value = 0.5086353321180399
$ws.Cells.Item(10, 17).Value = 540.0877496978508
$ws.Cells.Item(10, 18).Value = 4860.789747280656
$ws.Cells.Item(10, 19).Value = 0.1960154085704404
$ws.Cells.Item(10, 20).Value = 0.1960154085704404

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = 'MuSCs'
$ws.Cells.Item(11, 2).Value = 'Fgf1'
$ws.Cells.Item(11, 3).Value = 'Nrp1'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 4.037194
$ws.Cells.Item(11, 8).Value = 12.111582
$ws.Cells.Item(11, 9).Value = 0.3853751326204581
$ws.Cells.Item(11, 10).Value = 0.3853751326204581
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 52.441971
$ws.Cells.Item(11, 14).Value = 157.325913
$ws.Cells.Item(11, 15).Value = 0.1993888292903622
$ws.Cells.Item(11, 16).Value = 0.1993888292903622
$ws.Cells.Item(11, 17).Value = 211.718410669374
$ws.Cells.Item(11, 18).Value = 1905.465696024366
$ws.Cells.Item(11, 19).Value = 0.07683949653081122
$ws.Cells.Item(11, 20).Value = 0.07683949653081121

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 1).Value = 'MuSCs'
$ws.Cells.Item(12, 2).Value = 'Fgf1'
$ws.Cells.Item(12, 3).Value = 'Nrp1'
$ws.Cells.Item(12, 4).Value = 'MuSCs'
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 4.037194
$ws.Cells.Item(12, 8).Value = 12.111582
$ws.Cells.Item(12, 9).Value = 0.3853751326204581
$ws.Cells.Item(12, 10).Value = 0.3853751326204581
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 21.197691
$ws.Cells.Item(12, 14).Value = 63.593073
$ws.Cells.Item(12, 15).Value = 0.08059542216956049
$ws.Cells.Item(12, 16).Value = 0.08059542216956046
$ws.Cells.Item(12, 17).Value = 85.57919091905401
$ws.Cells.Item(12, 18).Value = 770.2127182714861
$ws.Cells.Item(12, 19).Value = 0.03105947150719618
$ws.Cells.Item(12, 20).Value = 0.03105947150719617

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = 'MuSCs'
$ws.Cells.Item(13, 2).Value = 'Fgf1'
$ws.Cells.Item(13, 3).Value = 'Nrp1'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 4.037194
$ws.Cells.Item(13, 8).Value = 12.111582
$ws.Cells.Item(13, 9).Value = 0.3853751326204581
$ws.Cells.Item(13, 10).Value = 0.3853751326204581
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 55.59592133333333
$ws.Cells.Item(13, 14).Value = 166.787764
$ws.Cells.Item(13, 15).Value = 0.2113804164220374
$ws.Cells.Item(13, 16).Value = 0.2113804164220373
$ws.Cells.Item(13, 17).Value = 224.4515200314054
$ws.Cells.Item(13, 18).Value = 2020.063680282648
$ws.Cells.Item(13, 19).Value = 0.0814607560120103
$ws.Cells.Item(13, 20).Value = 0.08146075601201028

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(14, 2).Value = 'Fgf1'
$ws.Cells.Item(14, 3).Value = 'Nrp1'
$ws.Cells.Item(14, 4).Value = 'ECs'
$ws.Cells.Item(14, 5).Value = 1.0
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.08000233333333333
$ws.Cells.Item(14, 8).Value = 0.240007
$ws.Cells.Item(14, 9).Value = 0.007636717437477471
$ws.Cells.Item(14, 10).Value = 0.007636717437477472
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 133.7780026666667
$ws.Cells.Item(14, 14).Value = 401.334008
$ws.Cells.Item(14, 15).Value = 0.50863533211804
$ws.Cells.Item(14, 16).Value = 0.5086353321180399
$ws.Cells.Item(14, 17).Value = 10.70255236200622
$ws.Cells.Item(14, 18).Value = 96.32297125805601
$ws.Cells.Item(14, 19).Value = 0.003884304310102981
$ws.Cells.Item(14, 20).Value = 0.003884304310102981

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(15, 2).Value = 'Fgf1'
$ws.Cells.Item(15, 3).Value = 'Nrp1'
$ws.Cells.Item(15, 4).Value = 'FAPs'
$ws.Cells.Item(15, 5).Value = 1.0
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.08000233333333333
$ws.Cells.Item(15, 8).Value = 0.240007
$ws.Cells.Item(15, 9).Value = 0.007636717437477471
$ws.Cells.Item(15, 10).Value = 0.007636717437477472
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 52.441971
$ws.Cells.Item(15, 14).Value = 157.325913
$ws.Cells.Item(15, 15).Value = 0.1993888292903622
$ws.Cells.Item(15, 16).Value = 0.1993888292903622
$ws.Cells.Item(15, 17).Value = 4.195480044599
$ws.Cells.Item(15, 18).Value = 37.759320401391
$ws.Cells.Item(15, 19).Value = 0.001522676149479928
$ws.Cells.Item(15, 20).Value = 0.001522676149479928

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(16, 2).Value = 'Fgf1'
$ws.Cells.Item(16, 3).Value = 'Nrp1'
$ws.Cells.Item(16, 4).Value = 'MuSCs'
$ws.Cells.Item(16, 5).Value = 1.0
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.08000233333333333
$ws.Cells.Item(16, 8).Value = 0.240007
$ws.Cells.Item(16, 9).Value = 0.007636717437477471
$ws.Cells.Item(16, 10).Value = 0.007636717437477472
$ws.Cells.Item(16, 11).Value = 3.0
$ws.Cells.Item(16, 12).Value = 1.0
$ws.Cells.Item(16, 13).Value = 21.197691
$ws.Cells.Item(16, 14).Value = 63.593073
$ws.Cells.Item(16, 15).Value = 0.08059542216956049
$ws.Cells.Item(16, 16).Value = 0.08059542216956046
$ws.Cells.Item(16, 17).Value = 1.695864741279
$ws.Cells.Item(16, 18).Value = 15.262782671511
$ws.Cells.Item(16, 19).Value = 0.0006154844658631409
$ws.Cells.Item(16, 20).Value = 0.0006154844658631408

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(17, 2).Value = 'Fgf1'
$ws.Cells.Item(17, 3).Value = 'Nrp1'
$ws.Cells.Item(17, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(17, 5).Value = 1.0
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.08000233333333333
$ws.Cells.Item(17, 8).Value = 0.240007
$ws.Cells.Item(17, 9).Value = 0.007636717437477471
$ws.Cells.Item(17, 10).Value = 0.007636717437477472
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 12).Value = 1.0
$ws.Cells.Item(17, 13).Value = 55.59592133333333
$ws.Cells.Item(17, 14).Value = 166.787764
$ws.Cells.Item(17, 15).Value = 0.2113804164220374
$ws.Cells.Item(17, 16).Value = 0.2113804164220373
$ws.Cells.Item(17, 17).Value = 4.447803430483111
$ws.Cells.Item(17, 18).Value = 40.03023087434799
$ws.Cells.Item(17, 19).Value = 0.001614252512031422
$ws.Cells.Item(17, 20).Value = 0.001614252512031422

Write-Host "Final UsedRange:" ($ws.UsedRange.Address())
